# Updates the cryptos price/volume table (Sheet1) to the latest scraped
# snapshot: refreshed Price (col D) / Volume(1h) (col E) figures for most
# rows, plus two rank swaps (Cosmos<->Dai at rows 24/25, and
# VeChain<->Celestia at rows 38/39) where the coin, its link, price and
# volume all moved together.
#
# Price-column values are written with a leading apostrophe so Excel keeps
# them as literal text instead of re-parsing multi-dot "thousands" prices
# (e.g. "42.127.13") or losing significant trailing zeros / switching to
# scientific notation on plain-looking decimals (e.g. "1.00" -> 1,
# "14.30" -> 14.3, "0.0302" -> 3.02E-02). The final Style reset clears the
# quote-prefix styling so the cells end up unstyled, matching the rest of
# the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.127.13"
$ws.Range("E2").Value = "  -1.74%  "

$ws.Range("D3").Value = "'2.245.38"
$ws.Range("E3").Value = "  -1.75%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'247.98"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("E6").Value = "  -1.20%  "

$ws.Range("D7").Value = "'76.71"
$ws.Range("E7").Value = "  +4.17%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  -3.59%  "

$ws.Range("D10").Value = "'41.28"
$ws.Range("E10").Value = "  +5.19%  "

$ws.Range("E11").Value = "  -2.23%  "

$ws.Range("D12").Value = "'7.21"
$ws.Range("E12").Value = "  -2.83%  "

$ws.Range("E13").Value = "  -3.11%  "

$ws.Range("D14").Value = "'2.581.48"
$ws.Range("E14").Value = "  -1.87%  "

$ws.Range("D15").Value = "'14.85"
$ws.Range("E15").Value = "  -2.82%  "

$ws.Range("D16").Value = "'0.859"
$ws.Range("E16").Value = "  -1.64%  "

$ws.Range("D17").Value = "'2.246.30"
$ws.Range("E17").Value = "  -1.92%  "

$ws.Range("D18").Value = "'42.015.53"
$ws.Range("E18").Value = "  -1.87%  "

$ws.Range("D19").Value = "'0.0₃0985"
$ws.Range("E19").Value = "  -1.84%  "

$ws.Range("E20").Value = "  -2.80%  "

$ws.Range("D21").Value = "'71.91"

$ws.Range("D22").Value = "'2.31"
$ws.Range("E22").Value = "  +4.09%  "

$ws.Range("D23").Value = "'231.88"
$ws.Range("E23").Value = "  -2.12%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'11.33"
$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("D26").Value = "'3.67"
$ws.Range("E26").Value = "  -5.75%  "

$ws.Range("E27").Value = "  -4.82%  "

$ws.Range("D28").Value = "'7.24"
$ws.Range("E28").Value = "  +11.74%  "

$ws.Range("E29").Value = "  -1.51%  "

$ws.Range("D30").Value = "'169.17"
$ws.Range("E30").Value = "  +1.46%  "

$ws.Range("D31").Value = "'20.57"
$ws.Range("E31").Value = "  -2.21%  "

$ws.Range("D32").Value = "'33.81"
$ws.Range("E32").Value = "  +8.10%  "

$ws.Range("D33").Value = "'0.0832"
$ws.Range("E33").Value = "  +0.89%  "

$ws.Range("E34").Value = "  -4.96%  "

$ws.Range("D35").Value = "'0.126"
$ws.Range("E35").Value = "  -0.76%  "

$ws.Range("D36").Value = "'4.56"
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("D37").Value = "'4.92"
$ws.Range("E37").Value = "  +3.30%  "

$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "'14.30"
$ws.Range("E38").Value = "  -1.30%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.0302"
$ws.Range("E39").Value = "  -2.56%  "

$ws.Range("D40").Value = "'5.95"
$ws.Range("E40").Value = "  +0.29%  "

$ws.Range("D41").Value = "'2.19"
$ws.Range("E41").Value = "  -6.33%  "

$ws.Range("D42").Value = "'113.45"
$ws.Range("E42").Value = "  +13.78%  "

$ws.Range("E43").Value = "  -7.27%  "

$ws.Range("D44").Value = "'61.17"
$ws.Range("E44").Value = "  -1.14%  "

$ws.Range("E45").Value = "  -4.02%  "

$ws.Range("E46").Value = "  -2.78%  "

$ws.Range("D47").Value = "'0.998"
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("E48").Value = "  -2.94%  "

$ws.Range("E50").Value = "  -12.96%  "

$ws.Range("E51").Value = "  -1.24%  "

# Clear the quote-prefix formatting picked up from the apostrophe-prefixed
# text assignments above so the cells stay unstyled (no explicit `s`),
# same as every other data cell on the sheet.
$ws.Range("D2:E51").Style = "Normal"
